$wb = $excel.ActiveWorkbook

# "Architect. Design Phase Defects" is the second worksheet in this workbook.
$ws = $wb.Worksheets.Item(2)

# Clear the leftover formatting across the checklist block (header row 9
# through the last checklist row 26) before filling it in with the actual
# architectural-design review data.
$ws.Range("B9:E26").Style = "Normal"

# Checklist rows (Crt. No. / Checked Item / Doc. page-line / Comments).
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "A01"
$ws.Range("E10").Value = "Organizararea programului este clara, exista pachete separate pentru fiecare strat al arhitecturii"

$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "A02"
$ws.Range("E11").Value = "Partitionarea si layering-ul este corect"

$ws.Range("B12").Value = 3
$ws.Range("C12").Value = "A03"
$ws.Range("E12").Value = "Arhitectura permite realizarea tuturor cerintelor"

$ws.Range("B13").Value = 4
$ws.Range("C13").Value = "A04"
$ws.Range("D13").Value = "PizzaService"
$ws.Range("E13").Value = "La nivelul proiectului exista un singur service responabil de business layer si acesta incorporeaza toate subsistemele"

$ws.Range("B14").Value = 5
$ws.Range("C14").Value = "A05"
$ws.Range("E14").Value = "Nu exista la nivelul aplicatiei o strategie de gestiune a erorilor, nu exista clase separate pentru un anumit tip de eroare. Sunt tratate doar erorile IO"

$ws.Range("B15").Value = 6
$ws.Range("C15").Value = "A06"
$ws.Range("E15").Value = "MVC model este incorporat in proiect"

$ws.Range("B16").Value = 7
$ws.Range("C16").Value = "A07"
$ws.Range("D16").Value = "PizzaService"
$ws.Range("E16").Value = "Se acupa atat cu plata cat si cu afisarea meniului. Din numele clasei nu iti dai seama care este scopul serviciului, o denumire mai buna ar fi PizzaManagementOrdersService"

$ws.Range("B17").Value = 8
$ws.Range("C17").Value = "A08"
$ws.Range("E17").Value = "Exista descriere la clase ?? Nu cred"

$ws.Range("B18").Value = 9
$ws.Range("C18").Value = "A09"
$ws.Range("E18").Value = "Relatiile 1 to many nu sunt bine scrise intre PaymentRepository si Payment, intre MenuGUIController si OrdersGUI ar trebui sa fie 1:10, "

$ws.Range("B19").Value = 10
$ws.Range("C19").Value = "A10"
$ws.Range("E19").Value = "The key entity classes are consistent with business and model layers"

# Reviewer info for this checklist (Student 1 / effort hours).
$ws.Range("I3").Value = "Silvia Pirlea"
$ws.Range("J3").Value = 235

# Make this the active sheet / selected cell, as it was left after editing.
$ws.Activate() | Out-Null
$ws.Range("J3").Select() | Out-Null
